$wb = $excel.ActiveWorkbook

# Rename sheets to add " forecasting" suffix
$wsUni = $wb.Worksheets.Item("univariate")
$wsUni.Name = "univariate forecasting"

$wsMulti = $wb.Worksheets.Item("multivariate")
$wsMulti.Name = "multivariate forecasting"

# Update selected cell on the univariate forecasting sheet
$wsUni.Activate()
$wsUni.Range("D3").Select()

# Update selected cell on the multivariate forecasting sheet
$wsMulti.Activate()
$wsMulti.Range("I12").Select()
